# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (col E) list is re-sorted from descending (2009 -> 1805)
# to ascending (1805 -> 2009) order, and the "Valor Mora" (col F) figures are
# re-assigned to match the new row order (same multiset of values, shifted
# into new bands).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period labels for E16:E35 (were descending 2009..1805).
$periods = @(
    "1805","1807","1808","1902","1903","1904","1905","1906","1907","1909",
    "1910","1911","2001","2002","2004","2005","2006","2007","2008","2009"
)

# New "Valor Mora" figures for F16:F35, aligned to the reordered periods.
$valores = @(
    31249,31249,31249,31249,31249,31249,31249,31249,31249,
    33125,33125,33125,33125,33125,
    35112,35112,35112,35112,35112,
    32771
)

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
